$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-29 Thursday", "2024-03-01 Friday"),
    @("899÷3=", "820÷4="),
    @("162÷7=", "740÷9="),
    @("198÷4=", "176÷8="),
    @("432÷3=", "412÷7="),
    @("324÷5=", "644÷9="),
    @("842÷9=", "863÷8="),
    @("194÷9=", "795÷7="),
    @("308÷2=", "601÷6="),
    @("993÷3=", "546÷3="),
    @("988÷4=", "264÷8="),
    @("896÷7=", "534÷8="),
    @("964÷9=", "291÷3="),
    @("337÷4=", "935÷4="),
    @("340÷9=", "544÷6="),
    @("787÷7=", "278÷6="),
    @("602÷6=", "650÷9="),
    @("998÷2=", "345÷9="),
    @("128÷3=", "784÷8="),
    @("468÷5=", "513÷7="),
    @("867÷7=", "499÷3="),
    @("202÷8=", "951÷2="),
    @("728÷6=", "133÷8="),
    @("233÷6=", "683÷8="),
    @("661÷3=", "267÷9="),
    @("867÷9=", "623÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
